$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows (1-3), including the superseded "2023-07-18" /
# 45125-date entries.
$ws.Rows("1:3").Delete()

# Re-enter the surviving record ("2023-07-19" / 4785) down at row 5.
$dateCell = $ws.Range("A5")
$dateCell.NumberFormat = "@"
$dateCell.Value = "2023-07-19"
$dateCell.Style = "Normal"

$ws.Range("B5").Value = 4785
